$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the supplementary table 2 column headers (row 1, columns B:G)
$ws.Range("B1").Value = "Guava Gain - Culture Tests & NES FLP"
$ws.Range("C1").Value = "Guava Gain - LysoTracker CCS"
$ws.Range("D1").Value = "Guava Gain - LysoTracker NES"
$ws.Range("E1").Value = "CytPix Voltage - Culture Tests"
$ws.Range("F1").Value = "CytPix Wavelength (nm)"
$ws.Range("G1").Value = "Guava Wavelength (nm)"

# Restore the active cell selection to B2, matching the saved view state
$ws.Range("B2").Select()
